# Apply the "output generated at 456a3b4" update to 杭州-漫展信息.xlsx
#
# Summary of changes:
#  1. Sheet "展览" (sheet index 1): several "想去人数" (F column) view-count
#     values are refreshed for existing rows, a brand-new exhibition row is
#     inserted at row 28 (shifting the former rows 28-38 down to 29-39), and
#     the F column of the shifted rows / the final row are refreshed too.
#  2. Sheet "本地生活" (sheet index 3): two F-column view-counts refreshed.
#  3. Sheet "全部类型" (sheet index 4): the same view-counts refreshed in
#     place (this combined sheet does not get the new row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# --- 1a. Simple F-column (view count) refreshes for untouched rows ---
$ws1.Range("F3").Value  = 8067
$ws1.Range("F4").Value  = 1902
$ws1.Range("F5").Value  = 6486
$ws1.Range("F7").Value  = 2036
$ws1.Range("F8").Value  = 558
$ws1.Range("F10").Value = 18
$ws1.Range("F14").Value = 63
$ws1.Range("F15").Value = 8412
$ws1.Range("F20").Value = 1797
$ws1.Range("F25").Value = 7
$ws1.Range("F26").Value = 50

# --- 1b. Insert a brand new row at position 28, pushing the former ---
# --- rows 28-38 down to 29-39.                                    ---
$ws1.Rows.Item(28).Insert()

# Fix up the formatting of the newly inserted (blank) row so that it
# matches the rest of the table (Excel's row insert drops the border
# that column A normally carries).
$ws1.Range("A27:I27").Copy()
$ws1.Range("A28:I28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 1c. Fill in the new row 28 with the new exhibition's data ---
$ws1.Range("A28").Value = 27

# B28 looks like a date ("2024-10-04"); force it to be stored as plain
# text (matching how every other date cell in this sheet is stored)
# instead of letting Excel auto-convert it to a date serial number.
$ws1.Range("B28").NumberFormat = "@"
$ws1.Range("B28").Value = "2024-10-04"
$ws1.Range("B28").Style = "Normal"

$ws1.Range("C28").Value = "杭州·创世次元第五人格同人only展"
$ws1.Range("D28").Value = "小河路与桥弄街交叉口东北50米 桥西历史文化街区"
$ws1.Range("E28").Value = "2024.10.04 10:00-10.05 17:00"
$ws1.Range("F28").Value = 1
$ws1.Range("G28").Value = 75
$ws1.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=92141"
$ws1.Range("I28").Value = "//i1.hdslb.com/bfs/openplatform/202409/MMF3dkAw1725550270634.jpeg"

# --- 1d. Refresh the F column (view counts) of the rows that were ---
# --- shifted down by the insert, plus the newly exposed row 39.   ---
$ws1.Range("F29").Value = 12
$ws1.Range("F30").Value = 2019
$ws1.Range("F31").Value = 840
$ws1.Range("F32").Value = 460
$ws1.Range("F33").Value = 5
$ws1.Range("F34").Value = 10
$ws1.Range("F35").Value = 160
$ws1.Range("F36").Value = 140
$ws1.Range("F37").Value = 2
$ws1.Range("F38").Value = 95
$ws1.Range("F39").Value = 3959

# ---------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 2315
$ws3.Range("F3").Value = 705

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 2315
$ws4.Range("F3").Value  = 705
$ws4.Range("F6").Value  = 8067
$ws4.Range("F9").Value  = 1902
$ws4.Range("F10").Value = 6486
$ws4.Range("F11").Value = 2036
$ws4.Range("F13").Value = 558
$ws4.Range("F15").Value = 18
$ws4.Range("F22").Value = 63
$ws4.Range("F23").Value = 8412
$ws4.Range("F28").Value = 1797
$ws4.Range("F32").Value = 50
$ws4.Range("F34").Value = 12
$ws4.Range("F35").Value = 2019
$ws4.Range("F38").Value = 460
$ws4.Range("F41").Value = 140
$ws4.Range("F44").Value = 3959

$wb.Save()
